$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# (cellRef, newValue, forceText) triples scraped from the crypto-price-
# refresh diff. forceText=1 means the literal string looks like a plain
# number (e.g. "1.000", "316.64", "0.07462"); those cells must be
# switched to a Text format before the write so Excel keeps the exact
# original digits/trailing zeros instead of silently coercing the
# assignment to a real number (which would drop formatting like the
# trailing zero in "1.000" -> 1, or introduce float noise).
$updates = @(
    ,@('D2', '27.965.82', 0)
    ,@('E2', '  -2.46%  ', 0)
    ,@('D3', '1.794.28', 0)
    ,@('E3', '  -0.64%  ', 0)
    ,@('D4', '1.000', 1)
    ,@('E4', '  -0.05%  ', 0)
    ,@('D5', '316.64', 1)
    ,@('E5', '  -0.33%  ', 0)
    ,@('E6', '  +0.01%  ', 0)
    ,@('D7', '0.5314', 1)
    ,@('E7', '  -2.82%  ', 0)
    ,@('D8', '0.3931', 1)
    ,@('E8', '  +3.39%  ', 0)
    ,@('D9', '0.07462', 1)
    ,@('E9', '  -0.78%  ', 0)
    ,@('D10', '41.52', 1)
    ,@('D11', '1.086', 1)
    ,@('E11', '  -2.72%  ', 0)
    ,@('E12', '  +0.00%  ', 0)
    ,@('B13', 'Polkadot', 0)
    ,@('C13', 'https://coinranking.com/coin/25W7FG7om+polkadot-dot', 0)
    ,@('D13', '6.184', 1)
    ,@('E13', '  +0.35%  ', 0)
    ,@('B14', 'Chainlink', 0)
    ,@('C14', 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link', 0)
    ,@('D14', '7.474', 1)
    ,@('E14', '  +0.95%  ', 0)
    ,@('D15', '20.36', 1)
    ,@('E15', '  -1.71%  ', 0)
    ,@('D16', '1.792.48', 0)
    ,@('E16', '  -0.25%  ', 0)
    ,@('D17', '88.37', 1)
    ,@('E17', '  -2.05%  ', 0)
    ,@('D18', '0.00001060', 1)
    ,@('E18', '  -0.81%  ', 0)
    ,@('D19', '0.06583', 1)
    ,@('E19', '  +1.48%  ', 0)
    ,@('D20', '1.000', 1)
    ,@('E20', '  +0.01%  ', 0)
    ,@('D21', '17.20', 1)
    ,@('E21', '  -1.03%  ', 0)
    ,@('D22', '5.955', 1)
    ,@('E22', '  +0.10%  ', 0)
    ,@('D23', '27.969.79', 0)
    ,@('E23', '  -2.44%  ', 0)
    ,@('D24', '11.09', 1)
    ,@('E24', '  -0.40%  ', 0)
    ,@('D25', '2.090', 1)
    ,@('E25', '  -0.17%  ', 0)
    ,@('D26', '156.66', 1)
    ,@('E26', '  -2.38%  ', 0)
    ,@('D27', '20.19', 1)
    ,@('E27', '  -1.30%  ', 0)
    ,@('D28', '2.003.15', 0)
    ,@('E28', '  +0.23%  ', 0)
    ,@('D29', '2.295', 1)
    ,@('E29', '  -2.93%  ', 0)
    ,@('D30', '122.01', 1)
    ,@('E30', '  -1.04%  ', 0)
    ,@('D31', '0.1085', 1)
    ,@('E31', '  +2.64%  ', 0)
    ,@('E32', '  -2.36%  ', 0)
    ,@('D33', '3.675', 1)
    ,@('E33', '  -0.28%  ', 0)
    ,@('D34', '5.507', 1)
    ,@('E34', '  -2.55%  ', 0)
    ,@('D35', '0.07072', 1)
    ,@('E35', '  +5.84%  ', 0)
    ,@('D36', '0.2216', 1)
    ,@('E36', '  -2.08%  ', 0)
    ,@('D37', '5.109', 1)
    ,@('E37', '  +1.51%  ', 0)
    ,@('D38', '0.02272', 1)
    ,@('E38', '  -1.39%  ', 0)
    ,@('D39', '8.392', 1)
    ,@('E39', '  -5.14%  ', 0)
    ,@('D40', '11.24', 1)
    ,@('E40', '  -0.34%  ', 0)
    ,@('D41', '1.183', 1)
    ,@('E41', '  -1.17%  ', 0)
    ,@('D42', '0.6118', 1)
    ,@('E42', '  -2.22%  ', 0)
    ,@('D43', '1.416', 1)
    ,@('E43', '  -1.47%  ', 0)
    ,@('D44', '13.29', 1)
    ,@('E44', '  -0.05%  ', 0)
    ,@('D45', '3.675', 1)
    ,@('E45', '  -0.58%  ', 0)
    ,@('D46', '0.5715', 1)
    ,@('E46', '  -2.56%  ', 0)
    ,@('D47', '125.35', 1)
    ,@('E47', '  -1.21%  ', 0)
    ,@('D48', '1.182', 1)
    ,@('E48', '  +1.91%  ', 0)
    ,@('D49', '1.918', 1)
    ,@('E49', '  -1.71%  ', 0)
    ,@('D50', '0.06806', 1)
    ,@('E50', '  -1.32%  ', 0)
    ,@('E51', '  -1.64%  ', 0)
)

foreach ($u in $updates) {
    $cellRef = $u[0]
    $newValue = $u[1]
    $forceText = $u[2]
    $rng = $ws.Range($cellRef)
    if ($forceText -eq 1) {
        $rng.NumberFormat = "@"
        $rng.Value = $newValue
        $rng.ClearFormats()
    } else {
        $rng.Value = $newValue
    }
}
